# Kazakhstan Premier League - swap match-detail data between mis-ordered row pairs.
# For each pair of rows below, the data in columns B:AD (id, HomeTeam, AwayTeam,
# scores, odds, etc.) was swapped between the two rows, while the leading
# index column A (and the header-less nothing else) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(99, 100),
    @(119, 120),
    @(129, 130),
    @(136, 137),
    @(143, 144),
    @(148, 149)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
